$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at position 265 (pushes old rows 265-281 down to 268-284)
$ws.Range("A265:A267").EntireRow.Insert()

# New row 265: Naranja, Lane Late, Primera
$ws.Range("A265").Value = 4
$ws.Range("B265").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C265").Value = "Los Lagos"
$ws.Range("D265").Value = 44509
$ws.Range("E265").Value = 10
$ws.Range("F265").Value = "Fruta"
$ws.Range("G265").Value = 100102
$ws.Range("H265").Value = "Cítricos"
$ws.Range("I265").Value = 100102005
$ws.Range("J265").Value = "Naranja"
$ws.Range("K265").Value = "Lane Late"
$ws.Range("L265").Value = "Primera"
$ws.Range("M265").Value = 400
$ws.Range("N265").Value = 13000
$ws.Range("O265").Value = 14000
$ws.Range("P265").Value = 13500
$ws.Range("Q265").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R265").Value = "Región de O'Higgins"
$ws.Range("S265").Value = 900
$ws.Range("T265").Value = 15

# New row 266: Naranja, Lane Late, Segunda
$ws.Range("A266").Value = 4
$ws.Range("B266").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C266").Value = "Los Lagos"
$ws.Range("D266").Value = 44509
$ws.Range("E266").Value = 10
$ws.Range("F266").Value = "Fruta"
$ws.Range("G266").Value = 100102
$ws.Range("H266").Value = "Cítricos"
$ws.Range("I266").Value = 100102005
$ws.Range("J266").Value = "Naranja"
$ws.Range("K266").Value = "Lane Late"
$ws.Range("L266").Value = "Segunda"
$ws.Range("M266").Value = 600
$ws.Range("N266").Value = 7500
$ws.Range("O266").Value = 8000
$ws.Range("P266").Value = 7750
$ws.Range("Q266").Value = "`$/malla 18 kilos"
$ws.Range("R266").Value = "Región de O'Higgins"
$ws.Range("S266").Value = 431
$ws.Range("T266").Value = 18

# New row 267: Naranja, Navel Late, Primera
$ws.Range("A267").Value = 4
$ws.Range("B267").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C267").Value = "Los Lagos"
$ws.Range("D267").Value = 44509
$ws.Range("E267").Value = 10
$ws.Range("F267").Value = "Fruta"
$ws.Range("G267").Value = 100102
$ws.Range("H267").Value = "Cítricos"
$ws.Range("I267").Value = 100102005
$ws.Range("J267").Value = "Naranja"
$ws.Range("K267").Value = "Navel Late"
$ws.Range("L267").Value = "Primera"
$ws.Range("M267").Value = 600
$ws.Range("N267").Value = 13000
$ws.Range("O267").Value = 14000
$ws.Range("P267").Value = 13500
$ws.Range("Q267").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R267").Value = "Provincia de Melipilla"
$ws.Range("S267").Value = 900
$ws.Range("T267").Value = 15
